# 3.1.2 chapter edit: "forklaring output = out-put"
#
# 1) "Input Output bulker." -> split "Output" into "Out" / "-" / "put"
#    (hyphenated "Out-put"), i.e. three runs: "...Input Out", "-", "put bulker."
# 2) "A_Do8_Dokumentobjekter_uten_dokumentfiler.xq" -> bold the
#    "Dokumentobjekter" portion of the filename.

$d = $word.ActiveDocument

# --- Edit 1: "Input Output bulker." -> "Input Out-put bulker." -------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Input Output bulker.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    # "Input Output bulker." -> rng1.Start .. rng1.End
    # "Output" starts right after "Input " (6 + 1 chars in)
    $outputStart = $rng1.Start + 6
    $outStart = $outputStart
    $outEnd = $outputStart + 3        # "Out"
    $putStart = $outputStart + 3      # "put" begins where "Out" ends
    $putEnd = $outputStart + 6        # end of "Output"

    # Insert the literal hyphen between "Out" and "put"
    $putRange = $d.Range($putStart, $putEnd)
    $putRange.InsertBefore("-")

    # Force a genuine run split at the hyphen (otherwise identically
    # formatted neighbouring text gets silently re-merged into one run)
    # by touching and releasing a character formatting toggle on just
    # the inserted hyphen.
    $hyphenRange = $d.Range($putStart, $putStart + 1)
    $hyphenRange.Font.Bold = $true
    $hyphenRange.Font.Bold = $false
}

# --- Edit 2: bold "Dokumentobjekter" inside the .xq filename ---------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("A_Do8_Dokumentobjekter_uten_dokumentfiler.xq", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $fileStart = $rng2.Start
    $boldStart = $fileStart + "A_Do8_".Length
    $boldEnd = $boldStart + "Dokumentobjekter".Length
    $boldRange = $d.Range($boldStart, $boldEnd)
    $boldRange.Font.Bold = $true
}
